$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.210.68"
$ws.Range("E2").Value = "  +1.39%  "
$ws.Range("D3").Value = "1.645.95"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("E6").Value = "  +2.06%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").Value = "  +1.44%  "
$ws.Range("E9").Value = "  +1.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.23%  "
$ws.Range("D12").Value = "1.877.33"
$ws.Range("D13").Value = "1.636.96"
$ws.Range("E14").Value = "  +0.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.542"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.83%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.69"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.28%  "
$ws.Range("D17").Value = "27.189.86"
$ws.Range("E17").Value = "  +1.24%  "
$ws.Range("E18").Value = "  +1.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "219.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.77%  "
$ws.Range("E21").Value = "  +5.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.84"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.18%  "
$ws.Range("E23").Value = "  +0.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.58"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.97%  "
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("E29").Value = "  -0.41%  "
$ws.Range("E30").Value = "  -0.57%  "
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("E32").Value = "  +0.57%  "
$ws.Range("E33").Value = "  +1.27%  "
$ws.Range("E34").Value = "  +1.91%  "
$ws.Range("D35").Value = "1.263.84"
$ws.Range("E35").Value = "  +1.23%  "
$ws.Range("E36").Value = "  +0.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0178"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.99%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.546"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.91%  "
$ws.Range("E39").Value = "  +1.66%  "
$ws.Range("E40").Value = "  -0.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.809"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.25"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.77%  "
$ws.Range("D44").Value = "1.786.95"
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.91"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.67"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("E47").Value = "  +0.85%  "
$ws.Range("E48").Value = "  +2.09%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.68"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.80%  "
$ws.Range("E51").Value = "  +0.22%  "
